$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.740.49'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '3.797.78'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('D4').Formula = "'0.999"
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Formula = "'702.63"
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Formula = "'169.54"
$ws.Range('E6').Value = '  -2.51%  '
$ws.Range('D7').Value = '3.795.52'
$ws.Range('E7').Value = '  -1.84%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Formula = "'0.522"
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('D11').Formula = "'7.59"
$ws.Range('E11').Value = '  +6.11%  '
$ws.Range('D12').Formula = "'0.458"
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('E13').Value = '  -3.51%  '
$ws.Range('D14').Formula = "'35.66"
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Value = '4.435.34'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '3.781.27'
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').Value = '70.694.39'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').Formula = "'17.31"
$ws.Range('E19').Value = '  -2.30%  '
$ws.Range('D20').Formula = "'7.09"
$ws.Range('E20').Value = '  -2.25%  '
$ws.Range('D21').Formula = "'494.91"
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').Formula = "'10.66"
$ws.Range('E22').Value = '  -4.59%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Formula = "'84.13"
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('E25').Value = '  -4.10%  '
$ws.Range('D26').Value = '3.947.08'
$ws.Range('E26').Value = '  -1.64%  '
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('D28').Formula = "'10.25"
$ws.Range('E28').Value = '  -5.52%  '
$ws.Range('D29').Formula = "'1.00"
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -7.08%  '
$ws.Range('D31').Formula = "'3.01"
$ws.Range('E31').Value = '  -6.12%  '
$ws.Range('D32').Formula = "'2.27"
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').Formula = "'7.30"
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('D34').Formula = "'28.94"
$ws.Range('E34').Value = '  -2.76%  '
$ws.Range('D35').Formula = "'0.176"
$ws.Range('E35').Value = '  -3.24%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '3.768.69'
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('D38').Formula = "'9.03"
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('D40').Formula = "'2.36"
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').Formula = "'5.93"
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('E44').Value = '  -6.28%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Formula = "'166.67"
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('E47').Value = '  +0.83%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Formula = "'420.90"
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').Formula = "'8.56"
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('E51').Value = '  -4.02%  '
